$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells retain their original Text type (matches source inlineStr cells);
# otherwise Excel auto-coerces numeric-looking strings (e.g. "1.00" -> 1) and
# strips significant trailing zeros / formatting.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '52.282.21'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.27%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.838.89'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.48%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '362.13'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.28%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '112.49'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.58%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.570'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +4.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.605'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.73%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.09'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.07%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0870'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.11%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.01%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.07'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.67%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.83'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.27%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.289.21'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.66%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.798.92'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.10%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.937'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +5.25%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '52.196.33'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.19%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.57'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.24%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.58%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.47'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0999'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.03%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '272.88'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.09%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.59'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.70%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.12%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.01'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.85%  '

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.06%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.35'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.94%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.18%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.01%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0482'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +8.66%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '35.37'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.96%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '52.50'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.65%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.90'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.51%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +13.17%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0852'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.48%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.11%  '

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.36%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.60%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.48'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.13%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.117'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.68%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.56'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.87%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '126.95'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.45%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '23.07'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.57%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.53%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.090.02'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.74%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.35'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.31'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.76%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.89'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.88%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.969'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.98%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.25'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.12%  '
